$wb = $excel.ActiveWorkbook

# -- "Image ID" sheet: add new row for "speed potion" (id 3) --
$wsImage = $wb.Worksheets.Item("Image ID")
$wsImage.Range("A5").Value = 3
$wsImage.Range("B5").Value = "speed potion"

# widen column B to fit the new text, and select the whole column
$wsImage.Columns("B").ColumnWidth = 12 + 1/6
$wsImage.Columns("B").Select()

# -- "Object Code" sheet: add new row for "PowerUp" (code 10) --
$wsCode = $wb.Worksheets.Item("Object Code")
$wsCode.Range("A6").Value = 10
$wsCode.Range("B6").Value = "PowerUp"

# leave the "Object Code" sheet as the active tab with B6 selected
$wsCode.Range("B6").Select()
